$wb = $excel.ActiveWorkbook

# --- Sheet 1: Metadata ---
$meta = $wb.Worksheets.Item("Metadata")

# Title (row 5)
$meta.Range("B5").Value = "Fetal Presentations"

# Date (row 8)
$meta.Range("B8").Value = "2024-02-27T09:44:15-05:00"

# Description (row 12) - now spans multiple lines and adds a sentence about NCHS
$desc = "This valueset contains codes to represent fetal presentations. This valueset is based on `n[PHVS_FetalPresentations_NCHS](https://phinvads.cdc.gov/vads/ViewValueSet.action?id=3C696B7B-BB33-4818-8996-1E3461E3F512).  `nOnly these codes are supported for exchange with NCHS."
$meta.Range("B12").Value = $desc

# --- Sheet 2: Include from SNOMED CT ---
# Original rows: 1 header, 2=70028003/Vertex, 3=6096002/Breech,
#                4=163518000/unsure-finding, 5=empty, 6=System URI/snomed url
# Target rows:   1 header, 2=70028003/Vertex, 3=6096002/Breech,
#                4=empty, 5=System URI/snomed url   (row 6 removed)
$snomed = $wb.Worksheets.Item("Include from SNOMED CT")

# Move the "System URI" row (row 6) up into row 5, preserving its formatting
$snomed.Range("A6").Copy($snomed.Range("A5"))
$snomed.Range("B6").Copy($snomed.Range("B5"))

# Clear row 4 (the removed "163518000" concept)
$snomed.Range("A4").Value = ""
$snomed.Range("B4").Value = ""

# Delete the now-duplicated trailing row 6
$snomed.Rows.Item(6).Delete()

# --- Sheet 3: Include from NullFlavor ---
# Original rows: 1 header, 2=OTH/Other, 3=empty, 4=System URI/NullFlavor url
# Target rows:   1 header, 2=OTH/Other, 3=UNK/unknown, 4=empty, 5=System URI/NullFlavor url
$nullFlavor = $wb.Worksheets.Item("Include from NullFlavor")

# Move the "System URI" row (row 4) down into row 5, preserving its formatting
$nullFlavor.Range("A4").Copy($nullFlavor.Range("A5"))
$nullFlavor.Range("B4").Copy($nullFlavor.Range("B5"))

# Clear row 4 so the new "UNK"/"unknown" row can take row 3 while row 4 stays the blank row
$nullFlavor.Range("A4").Value = ""
$nullFlavor.Range("B4").Value = ""

# Insert the new "UNK" / "unknown" concept into row 3
$nullFlavor.Range("A3").Value = "UNK"
$nullFlavor.Range("B3").Value = "unknown"
